$wb = $excel.ActiveWorkbook

# --- TestData sheet: move selection from D17 to C25 ---
$wsTestData = $wb.Worksheets.Item("TestData")
$wsTestData.Activate()
$wsTestData.Range("C25").Select()

# --- NewReportTestData sheet: change I61 from Yes to No ---
$wsNewReport = $wb.Worksheets.Item("NewReportTestData")
$wsNewReport.Range("I61").Value = "No"

# --- DoubleClickReportTestData sheet: insert new column L (ApplyFeesfromcodes)
#     and append new column S (TestCase) ---
$wsDC = $wb.Worksheets.Item("DoubleClickReportTestData")

# Category column (C) text change -- set first so the new shared string
# "DoubleClick Reports" is registered before "ApplyFeesfromcodes"/"TestCase"
$wsDC.Range("C2").Value = "DoubleClick Reports"
$wsDC.Range("C3").Value = "DoubleClick Reports"

# Insert a new column before L (shifts L..Q to M..R)
$wsDC.Columns.Item(12).Insert()

# New column L header + values
$wsDC.Range("L1").Value = "ApplyFeesfromcodes"
$wsDC.Range("L2").Value = "Yes"
$wsDC.Range("L3").Value = "No"

# ExecutionStatus (now column R after shift) changes from Yes to No
$wsDC.Range("R2").Value = "No"
$wsDC.Range("R3").Value = "No"

# New column S (TestCase) header + values
$wsDC.Range("S1").Value = "TestCase"
$wsDC.Range("S1").Font.Bold = $true
$wsDC.Range("S1").NumberFormat = "@"
$wsDC.Range("S2").Value = "Test1"
$wsDC.Range("S3").Value = "Test2"

# Column B best-fit width (date column)
$wsDC.Columns.Item(2).AutoFit()

# View state: selection on DoubleClick sheet
$wsDC.Range("M2").Select()

# Make NewReportTestData the active sheet/tab, restore its pane/selection
$wsNewReport.Activate()
$wsNewReport.Range("J3").Select()
